$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header for the joined column
$ws.Range("D1").Value = "sample_type"

# Fill D2:D57 with the constant "permanent" value from the join
$lastRow = 57
$ws.Range("D2:D$lastRow").Value = "permanent"

# Update the active selection/view to reflect the edit location (F41)
$ws.Range("F41").Select()
